# Auto-generated Excel COM-interop script (PowerShell-style)
# Implements the "New crime data collected" edit to the 72nd Precinct
# CompStat worksheet: bump the report Volume/Number and the week-covering
# date range, and refresh every weekly/28-day/YTD/2yr crime-stat cell in
# rows 15-30 to the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header: volume/number and the week-covering date range ---
$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Cells that change "flavor" (numeric <-> text placeholder) ---
# These are first cloned (value + style) from an untouched donor cell of
# the desired flavor/style elsewhere on the sheet, which is how Excel COM
# reuses an existing style index / shared-string entry instead of minting a
# new one; numeric targets then get their real value poked in afterwards.
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("F15").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("L14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("F15").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 2
$ws.Range("L14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("F15").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("C14").Copy($ws.Range("G29"))
$ws.Range("E14").Copy($ws.Range("H29"))
$ws.Range("C14").Copy($ws.Range("G30"))
$ws.Range("E14").Copy($ws.Range("H30"))

# --- Plain value updates (same cell flavor/style as before) ---
$ws.Range("L15").Value = -83.333333333333
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 6
$ws.Range("H16").Value = -57.142857142857
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 28
$ws.Range("K16").Value = -53.571428571428
$ws.Range("L16").Value = -45.833333333333
$ws.Range("M16").Value = -62.857142857142
$ws.Range("N16").Value = -94.170403587443
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -90.90909090909
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = -26.923076923076
$ws.Range("I17").Value = 38
$ws.Range("J17").Value = 51
$ws.Range("K17").Value = -25.490196078431
$ws.Range("L17").Value = 22.58064516129
$ws.Range("M17").Value = 46.153846153846
$ws.Range("N17").Value = -48.648648648648
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -80
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = -67.647058823529
$ws.Range("L18").Value = -60.714285714285
$ws.Range("M18").Value = -75.555555555555
$ws.Range("N18").Value = -93.67816091954
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 15
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -45.098039215686
$ws.Range("I19").Value = 43
$ws.Range("J19").Value = 99
$ws.Range("K19").Value = -56.565656565656
$ws.Range("L19").Value = -52.222222222222
$ws.Range("M19").Value = -14
$ws.Range("N19").Value = -30.645161290322
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -72.727272727272
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 31
$ws.Range("K20").Value = -77.419354838709
$ws.Range("L20").Value = -69.565217391304
$ws.Range("M20").Value = -70.833333333333
$ws.Range("N20").Value = -96.551724137931
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -79.487179487179
$ws.Range("F21").Value = 60
$ws.Range("G21").Value = 119
$ws.Range("H21").Value = -49.579831932773
$ws.Range("I21").Value = 113
$ws.Range("J21").Value = 246
$ws.Range("K21").Value = -54.065040650406
$ws.Range("L21").Value = -44.334975369458
$ws.Range("M21").Value = -38.586956521739
$ws.Range("N21").Value = -84.811827956989
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = -66.666666666666
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = -16.666666666666
$ws.Range("I24").Value = 118
$ws.Range("J24").Value = 158
$ws.Range("K24").Value = -25.316455696202
$ws.Range("L24").Value = -32.571428571428
$ws.Range("M24").Value = -15.714285714285
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -45.454545454545
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = -30
$ws.Range("I25").Value = 39
$ws.Range("J25").Value = 78
$ws.Range("K25").Value = -50
$ws.Range("L25").Value = -62.135922330097
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -23.076923076923
$ws.Range("F26").Value = 47
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = 42.424242424242
$ws.Range("I26").Value = 79
$ws.Range("J26").Value = 76
$ws.Range("K26").Value = 3.947368421052
$ws.Range("L26").Value = 9.722222222222
$ws.Range("M26").Value = -17.708333333333
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = -20
$ws.Range("L27").Value = -42.857142857142
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -80
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = 14.285714285714
$ws.Range("L28").Value = -20

